# Applies the "Bugfixes. Added some standard-functionality." commit:
#  - TODO CMS sheet: mark the "Logout" row as done, and append a
#    new open TODO row for searching an order by its orderID.

$wb = $excel.ActiveWorkbook

$wsCms = $wb.Worksheets.Item("TODO CMS")

# Row 3 ("Logout") switches Status from "offen" to "done".
# Copy the formatting from an existing "done" cell (B9) so the green
# status fill/font come across exactly like Excel's own formatting.
$wsCms.Range("B9").Copy() | Out-Null
$wsCms.Range("B3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$wsCms.Range("B3").Value = "done"

# New row 11: a fresh open TODO item
$wsCms.Range("A11").Value = "Order anhand orderID suchen können"
$wsCms.Range("B11").Value = "offen"
$wsCms.Range("B11").Style = "Schlecht"
$wsCms.Range("C11").Value = "Jonas"

$wsCms.Range("D10").Select()
